$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix style on C85:C89 (apply right-align style, matches existing column C style) ---
$ws.Range("C85:C89").HorizontalAlignment = -4152

# --- Write new data rows 90:117 (columns A, B, C) ---
$ws.Range("A90").Value = 45218.666805555556
$ws.Range("B90").Value = 0
$ws.Range("C90").Value = 166.8
$ws.Range("A91").Value = 45218.708472222221
$ws.Range("B91").Value = 0
$ws.Range("C91").Value = 166.7
$ws.Range("A92").Value = 45218.750138888892
$ws.Range("B92").Value = 0
$ws.Range("C92").Value = 166.7
$ws.Range("A93").Value = 45218.791805555556
$ws.Range("B93").Value = 0
$ws.Range("C93").Value = 166.4
$ws.Range("A94").Value = 45218.833472222221
$ws.Range("B94").Value = 0
$ws.Range("C94").Value = 166.9
$ws.Range("A95").Value = 45218.875138888892
$ws.Range("B95").Value = 0
$ws.Range("C95").Value = 166.3
$ws.Range("A96").Value = 45218.916805555556
$ws.Range("B96").Value = 0
$ws.Range("C96").Value = 166.4
$ws.Range("A97").Value = 45218.958472222221
$ws.Range("B97").Value = 0
$ws.Range("C97").Value = 166.7
$ws.Range("A98").Value = 45219.000138888892
$ws.Range("B98").Value = 0
$ws.Range("C98").Value = 166.6
$ws.Range("A99").Value = 45219.041805555556
$ws.Range("B99").Value = 0
$ws.Range("C99").Value = 166.8
$ws.Range("A100").Value = 45219.083472222221
$ws.Range("B100").Value = 0
$ws.Range("C100").Value = 166.7
$ws.Range("A101").Value = 45219.125138888892
$ws.Range("B101").Value = 0
$ws.Range("C101").Value = 166.7
$ws.Range("A102").Value = 45219.166805555556
$ws.Range("B102").Value = 0
$ws.Range("C102").Value = 166.4
$ws.Range("A103").Value = 45219.208472222221
$ws.Range("B103").Value = 0
$ws.Range("C103").Value = 166
$ws.Range("A104").Value = 45219.250138888892
$ws.Range("B104").Value = 0
$ws.Range("C104").Value = 165.5
$ws.Range("A105").Value = 45219.291805555556
$ws.Range("B105").Value = 0
$ws.Range("C105").Value = 165.2
$ws.Range("A106").Value = 45219.333472222221
$ws.Range("B106").Value = 0
$ws.Range("C106").Value = 165.1
$ws.Range("A107").Value = 45219.375138888892
$ws.Range("B107").Value = 0
$ws.Range("C107").Value = 164.8
$ws.Range("A108").Value = 45219.416805555556
$ws.Range("B108").Value = 0
$ws.Range("C108").Value = 164.6
$ws.Range("A109").Value = 45219.458472222221
$ws.Range("B109").Value = 0
$ws.Range("C109").Value = 164.5
$ws.Range("A110").Value = 45219.500138888892
$ws.Range("B110").Value = 0
$ws.Range("C110").Value = 164.3
$ws.Range("A111").Value = 45219.541805555556
$ws.Range("B111").Value = 0
$ws.Range("C111").Value = 164.2
$ws.Range("A112").Value = 45219.583472222221
$ws.Range("B112").Value = 0
$ws.Range("C112").Value = 163.9
$ws.Range("A113").Value = 45219.625138888892
$ws.Range("B113").Value = 0
$ws.Range("C113").Value = 164.4
$ws.Range("A114").Value = 45219.666805555556
$ws.Range("B114").Value = 0
$ws.Range("C114").Value = 164.2
$ws.Range("A115").Value = 45219.708472222221
$ws.Range("B115").Value = 0
$ws.Range("C115").Value = 164.1
$ws.Range("A116").Value = 45219.750138888892
$ws.Range("B116").Value = 0
$ws.Range("C116").Value = 164.6
$ws.Range("A117").Value = 45219.791805555556
$ws.Range("B117").Value = 1.2
$ws.Range("C117").Value = 164.2

# --- Copy formatting (number formats / alignment) from row 89 down to the new rows ---
$ws.Range("A89:C89").Copy()
$ws.Range("A90:C117").PasteSpecial(-4122)

# --- Add empty, styled D column cells for rows 98:117 (style matches column C) ---
$ws.Range("C1").Copy()
$ws.Range("D98:D117").PasteSpecial(-4122)

# --- Update the sheet view: selection + scroll position ---
$ws.Range("G93").Select()
$excel.ActiveWindow.ScrollRow = 100
$excel.ActiveWindow.ScrollColumn = 1
